# Scene 30A edit: tag Roxy's dialogue lines with her expression/pose cues,
# matching the pattern already used for Prim/Mara lines in this scene, and
# add a (blank) default header to the section.

$d = $word.ActiveDocument

function Insert-CueParagraphBefore($anchorText, $cueText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $idx = $rng.Paragraphs.Item(1).Index
    $rng.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.Text = $cueText
}

function Insert-CueParagraphAfter($anchorText, $cueText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $idx = $rng.Paragraphs.Item(1).Index
    $rng.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $cueText
}

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- Insert new stand-alone cue paragraphs (anchored on the still-unedited
#     original text, so these run before the wording replacements below). ---

Insert-CueParagraphAfter "Roxy: Remember me?" "Roxy (neutral curious):"
Insert-CueParagraphBefore "Pro: ...right?" "Roxy (neutral disappointed):"
Insert-CueParagraphBefore "Unexpectedly, she lets out a laugh." "Roxy (neutral laughing):"
Insert-CueParagraphBefore "She laughs again, but this time it sounds more uneasy." "Roxy (neutral laughing):"
Insert-CueParagraphBefore "She leaves the auditorium, stopping to talk with her friend for a bit." "Roxy (exit):"

# --- Tag Roxy's existing dialogue lines with her expression/pose. ---

Replace-Text "Roxy: Hey, you’re here again." "Roxy (neutral curious): Hey, you’re here again."
Replace-Text "Roxy: Remember me?" "Roxy (neutral smiling): Remember me?"
Replace-Text "Roxy: Not quite…" "Roxy (neutral smiling_nervous): Not quite…"
Replace-Text "Roxy: I’m just kidding. I’m glad you remembered who I am." "Roxy (neutral grinning): I’m just kidding. I’m glad you remembered who I am."
Replace-Text "Roxy: You came here with Prim again?" "Roxy (neutral smiling): You came here with Prim again?"
Replace-Text "Roxy: Today we’re doing the final audition for the performance next week. I mean, Prim’s pretty much guaranteed to get the piano part, but I’m still a little nervous." "Roxy (neutral smiling_worried): Today we’re doing the final audition for the performance next week. I mean, Prim’s pretty much guaranteed to get the piano part, but I’m still a little nervous."
Replace-Text "Roxy: You sure you should be saying that?" "Roxy (neutral curious): You sure you should be saying that?"
Replace-Text "Roxy: Really…?" "Roxy (neutral disappointed): Really…?"
Replace-Text "Roxy: I kinda feel bad for Prim." "Roxy (neutral smiling_nervous): I kinda feel bad for Prim."
Replace-Text "Roxy: Well anyways, I should probably go get some practice in before we start." "Roxy (neutral smiling): Well anyways, I should probably go get some practice in before we start."
Replace-Text "Roxy: I’ll see you later, Pro." "Roxy (waving smiling): I’ll see you later, Pro."

# --- Give the section a (blank) default header. ---

$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$hdr.Range.InsertParagraphAfter()

Write-Output "done"
